$wb = $excel.ActiveWorkbook

# The simulation-name formula lives on the "TRANSIENT" sheet; fall back to
# the active sheet if for some reason the named lookup isn't available.
$ws = $wb.Worksheets.Item("TRANSIENT")
if (-not $ws) { $ws = $wb.ActiveSheet }

# Rename the test case: the TEXTJOIN formula in E3 appends a literal tag to
# the simulation name ("...refined" -> "...case2").
$ws.Range("E3").Formula = '=_xlfn.TEXTJOIN("_",TRUE,A6,E6,A8,E8,[1]GRID!$A$4,[1]GRID!$E$4,"case2")'
